$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Update "Success" -> "Pass" and "No_Success" -> "Fail"
$ws.Range("E3").Value = "Pass"
$ws.Range("E4").Value = "Fail"

# Update "Sam" -> "Tom"
$ws.Range("D9").Value = "Tom"

# Activate sheet and select F8, scrolled to top-left of A1
$ws.Activate()
$ws.Range("A1").Select()
$ws.Range("F8").Select()
